$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 375.5
$ws.Range("I33").Value = 318.875
$ws.Range("K33").Value = 318.875
$ws.Range("M33").Value = -89.875
$ws.Range("H62").Value = 2997.25
$ws.Range("I62").Value = 1906.7
$ws.Range("J62").Value = 8450
$ws.Range("K62").Value = 1906.7
$ws.Range("L62").Value = 8450
$ws.Range("M62").Value = -1282.7
$ws.Range("N62").Value = -9698
$ws.Range("H65").Value = 2997.25
$ws.Range("I65").Value = 1906.7
$ws.Range("J65").Value = 8450
$ws.Range("K65").Value = 9533.5
$ws.Range("L65").Value = 42250
$ws.Range("M65").Value = -6413.5
$ws.Range("N65").Value = -48490
$ws.Range("H100").Value = 28573928
$ws.Range("I100").Value = 40002400
$ws.Range("J100").Value = 2750
$ws.Range("K100").Value = 40002400
$ws.Range("L100").Value = 2750
$ws.Range("M100").Value = -40001859
$ws.Range("N100").Value = -3832
$ws.Range("H105").Value = 39713.57
$ws.Range("J105").Value = 39713.57
$ws.Range("L105").Value = 39713.57
$ws.Range("N105").Value = -46701.57
$ws.Range("H112").Value = 1268.3729
$ws.Range("I112").Value = 525
$ws.Range("J112").Value = 1322.4364
$ws.Range("K112").Value = 1575
$ws.Range("L112").Value = 3967.3092
$ws.Range("M112").Value = -467
$ws.Range("N112").Value = -6183.3092
$ws.Range("H132").Value = 40006480
$ws.Range("I132").Value = 45461044
$ws.Range("J132").Value = 6332
$ws.Range("K132").Value = 136383132
$ws.Range("L132").Value = 18996
$ws.Range("M132").Value = -136380602
$ws.Range("N132").Value = -24056
$ws.Range("H137").Value = 1289004
$ws.Range("I137").Value = 1702109.1
$ws.Range("J137").Value = 3787.7778
$ws.Range("K137").Value = 5106327.300000001
$ws.Range("L137").Value = 11363.3334
$ws.Range("M137").Value = -5103777.300000001
$ws.Range("N137").Value = -16463.3334
$ws.Range("H138").Value = 2702.26
$ws.Range("I138").Value = 834.46155
$ws.Range("J138").Value = 2981.3562
$ws.Range("K138").Value = 2503.38465
$ws.Range("L138").Value = 8944.068600000001
$ws.Range("M138").Value = 2636.61535
$ws.Range("N138").Value = -19224.0686

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H93").Value = 24499.334
$ws.Range("J93").Value = 24499.334
$ws.Range("L93").Value = 24499.334
$ws.Range("N93").Value = -29491.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 20000
$ws.Range("J58").Value = 20000
$ws.Range("L58").Value = 20000
$ws.Range("N58").Value = -20588
$ws.Range("H80").Value = 356.07693
$ws.Range("I80").Value = 345
$ws.Range("J80").Value = 361
$ws.Range("K80").Value = 345
$ws.Range("L80").Value = 361
$ws.Range("M80").Value = 653
$ws.Range("N80").Value = -2357
$ws.Range("H83").Value = 356.07693
$ws.Range("I83").Value = 345
$ws.Range("J83").Value = 361
$ws.Range("K83").Value = 1725
$ws.Range("L83").Value = 1805
$ws.Range("M83").Value = 3267
$ws.Range("N83").Value = -11789
$ws.Range("H134").Value = 3367.3333
$ws.Range("I134").Value = 1898.8572
$ws.Range("K134").Value = 5696.571599999999
$ws.Range("M134").Value = -3161.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3994.5
$ws.Range("I31").Value = 1892.8889
$ws.Range("K31").Value = 1892.8889
$ws.Range("M31").Value = -1597.8889
$ws.Range("H34").Value = 3994.5
$ws.Range("I34").Value = 1892.8889
$ws.Range("K34").Value = 1892.8889
$ws.Range("M34").Value = -1690.8889
$ws.Range("H105").Value = 1897.2
$ws.Range("I105").Value = 1208.2858
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1208.2858
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 538.7141999999999
$ws.Range("N105").Value = -5994
$ws.Range("H115").Value = 37800
$ws.Range("J115").Value = 37800
$ws.Range("L115").Value = 37800
$ws.Range("N115").Value = -40150
$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -47258
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -32620
$ws.Range("H122").Value = 2616.6365
$ws.Range("I122").Value = 1372.875
$ws.Range("K122").Value = 4118.625
$ws.Range("M122").Value = -1668.625
$ws.Range("H137").Value = 33367
$ws.Range("J137").Value = 33367
$ws.Range("L137").Value = 33367
$ws.Range("N137").Value = -43567

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2529.818
$ws.Range("I3").Value = 1691
$ws.Range("J3").Value = 4766.6665
$ws.Range("K3").Value = 5073
$ws.Range("L3").Value = 14299.9995
$ws.Range("M3").Value = -4961
$ws.Range("N3").Value = -14523.9995
$ws.Range("H5").Value = 446181.6
$ws.Range("I5").Value = 460.86667
$ws.Range("J5").Value = 891902.3
$ws.Range("K5").Value = 1382.60001
$ws.Range("L5").Value = 2675706.9
$ws.Range("M5").Value = -1270.60001
$ws.Range("N5").Value = -2675930.9
$ws.Range("H122").Value = 3208.0698
$ws.Range("I122").Value = 1183.3334
$ws.Range("J122").Value = 3536.4055
$ws.Range("K122").Value = 10650.0006
$ws.Range("L122").Value = 31827.6495
$ws.Range("M122").Value = -8200.000599999999
$ws.Range("N122").Value = -36727.6495
$ws.Range("H131").Value = 817.0947
$ws.Range("J131").Value = 822.06384
$ws.Range("L131").Value = 2466.19152
$ws.Range("N131").Value = -12546.19152
$ws.Range("H135").Value = 446181.6
$ws.Range("I135").Value = 460.86667
$ws.Range("J135").Value = 891902.3
$ws.Range("K135").Value = 4147.80003
$ws.Range("L135").Value = 8027120.7
$ws.Range("M135").Value = -1612.80003
$ws.Range("N135").Value = -8032190.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5260.164
$ws.Range("I70").Value = 5011.881
$ws.Range("J70").Value = 5809
$ws.Range("K70").Value = 5011.881
$ws.Range("L70").Value = 5809
$ws.Range("M70").Value = -4741.881
$ws.Range("N70").Value = -6349
$ws.Range("H73").Value = 5260.164
$ws.Range("I73").Value = 5011.881
$ws.Range("J73").Value = 5809
$ws.Range("K73").Value = 5011.881
$ws.Range("L73").Value = 5809
$ws.Range("M73").Value = -4075.881
$ws.Range("N73").Value = -7681
$ws.Range("H80").Value = 35718144
$ws.Range("I80").Value = 62502250
$ws.Range("J80").Value = 6000
$ws.Range("K80").Value = 62502250
$ws.Range("L80").Value = 6000
$ws.Range("M80").Value = -62501252
$ws.Range("N80").Value = -7996
$ws.Range("H83").Value = 35718144
$ws.Range("I83").Value = 62502250
$ws.Range("J83").Value = 6000
$ws.Range("K83").Value = 312511250
$ws.Range("L83").Value = 30000
$ws.Range("M83").Value = -312506258
$ws.Range("N83").Value = -39984
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
$ws.Range("H132").Value = 3078.7407
$ws.Range("I132").Value = 2381.4
$ws.Range("J132").Value = 5071.143
$ws.Range("K132").Value = 7144.200000000001
$ws.Range("L132").Value = 15213.429
$ws.Range("M132").Value = -4614.200000000001
$ws.Range("N132").Value = -20273.429
$ws.Range("H137").Value = 40320
$ws.Range("J137").Value = 40320
$ws.Range("L137").Value = 40320
$ws.Range("N137").Value = -50520

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3568.818
$ws.Range("I40").Value = 2425.7
$ws.Range("K40").Value = 2425.7
$ws.Range("M40").Value = -2289.7
$ws.Range("H136").Value = 5342.3887
$ws.Range("I136").Value = 1409
$ws.Range("J136").Value = 7845.4546
$ws.Range("K136").Value = 4227
$ws.Range("L136").Value = 23536.3638
$ws.Range("M136").Value = -1677
$ws.Range("N136").Value = -28636.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 61030.582
$ws.Range("J46").Value = 61030.582
$ws.Range("L46").Value = 61030.582
$ws.Range("N46").Value = -61492.582
$ws.Range("H108").Value = 28500
$ws.Range("J108").Value = 28500
$ws.Range("L108").Value = 28500
$ws.Range("N108").Value = -36180
$ws.Range("H132").Value = 37052228
$ws.Range("I132").Value = 51100
$ws.Range("J132").Value = 47623976
$ws.Range("K132").Value = 153300
$ws.Range("L132").Value = 142871928
$ws.Range("M132").Value = -150770
$ws.Range("N132").Value = -142876988
$ws.Range("H134").Value = 61030.582
$ws.Range("J134").Value = 61030.582
$ws.Range("L134").Value = 183091.746
$ws.Range("N134").Value = -188161.746
$ws.Range("H136").Value = 17163.637
$ws.Range("I136").Value = 51497.5
$ws.Range("K136").Value = 154492.5
$ws.Range("M136").Value = -151942.5
$ws.Range("H140").Value = 30751.2
$ws.Range("J140").Value = 30751.2
$ws.Range("L140").Value = 30751.2
$ws.Range("N140").Value = -41111.2
